$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-21 01:38:23"

foreach ($row in 2..8) {
    $ws.Range("A$row").Value = $newTimestamp
}
